$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 248
$ws1.Range("F3").Value = 69
$ws1.Range("F5").Value = 6485
$ws1.Range("F6").Value = 5240
$ws1.Range("F7").Value = 440
$ws1.Range("F8").Value = 62
$ws1.Range("F11").Value = 224
$ws1.Range("F12").Value = 37

# Sheet "全部类型" (All types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 248
$ws4.Range("F3").Value = 69
$ws4.Range("F5").Value = 6485
$ws4.Range("F6").Value = 5240
$ws4.Range("F7").Value = 440
$ws4.Range("F8").Value = 62
$ws4.Range("F11").Value = 224
$ws4.Range("F14").Value = 37
